# Update the timesheet to reflect the new resource/week/data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resource name
$ws.Range("B3").Value = "Vineet Rajput"

# Week (client unchanged at B4)
$ws.Range("B5").Value = "7/23/2018-7/29/2018"

# Day-of-week header row
$ws.Range("E7").Value = "23-Jul"
$ws.Range("F7").Value = "24-Jul"
$ws.Range("G7").Value = "25-Jul"
$ws.Range("H7").Value = "26-Jul"
$ws.Range("I7").Value = "27-Jul"
$ws.Range("J7").Value = "28-Jul"
$ws.Range("K7").Value = "29-Jul"

# Activity row
$ws.Range("B9").Value = "Client Call ( sanchit ) "
$ws.Range("E9").Value = "1"
$ws.Range("F9").Value = "1"
$ws.Range("G9").Value = "1"
$ws.Range("H9").Value = "1"
$ws.Range("I9").Value = "1"
$ws.Range("J9").Value = "1"
$ws.Range("K9").Value = "1"

# Per-day total row
$ws.Range("E18").Value = "1"
$ws.Range("F18").Value = "1"
$ws.Range("G18").Value = "1"
$ws.Range("H18").Value = "1"
$ws.Range("I18").Value = "1"
$ws.Range("J18").Value = "1"
$ws.Range("K18").Value = "1"

# Weekly total
$ws.Range("D19").Value = "7"

# Offshore consultant / Omar colon contact info updated to emails
$ws.Range("B22").Value = "Offshore Consultant's Project Manager's Name :dhananjayKumar@gmail.com"
$ws.Range("B25").Value = "Omarcolon@gmail.com"
